$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header cells I1 ("I0") and J1 ("IF"), copying the formatting
# (bold font, border, centered alignment) from the existing header cell H1
# so that the new headers reuse the same style as the rest of the header row.
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Column I ("I0") is a constant column of 1s for every data row.
$ws.Range("I2:I39").Value = 1

# Column J ("IF") duplicates the values already present in column H ("IP").
$ws.Range("J2:J39").Value = $ws.Range("H2:H39").Value2

# Clear the clipboard marquee left over from the Copy operation above.
$excel.CutCopyMode = 0
